$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — rows 4, 5, 11
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 7982
$wsExpo.Range("F5").Value = 5830
$wsExpo.Range("F11").Value = 358

# Sheet "全部类型" (All types) — rows 4, 5, 14
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7982
$wsAll.Range("F5").Value = 5830
$wsAll.Range("F14").Value = 358
